# Updated cryptos list on Sat May 27 15:40:36 UTC 2023 with GitHub Actions
# Applies latest price/volume(1h) snapshot to the cryptos worksheet (columns D and E, rows 2-51).
# Text values that look numeric are apostrophe-prefixed to keep them stored as text (matching
# the original inline-string cell type), then the cell style is reset to Normal so no stray
# quote-prefix / number-format style is left applied to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.998.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "1.846.08"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'1.012"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'309.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4767"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.97%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3684"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07231"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.9320"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'19.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.17%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.07739"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "1.835.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.387"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.71%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'6.469"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'88.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000008666"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "27.035.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'14.57"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("E22").Value = "  +0.64%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'10.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'1.961"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'153.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'18.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'2.010"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'114.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'4.972"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.74%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'3.325"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.50%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'0.7439"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'4.507"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'2.690"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.26%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.53%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.01964"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'0.05267"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'2.972"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.5253"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.81%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'7.028"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.75%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'8.314"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.4736"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("E46").Value = "  +0.61%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'101.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.49%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'1.610"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'65.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.06072"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.8934"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.61%  "
$ws.Range("E51").Style = "Normal"
